# Update "想去人数" (interested-people count) values in column F
# across the "展览", "演出" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1025
$ws1.Range("F6").Value = 164
$ws1.Range("F8").Value = 194
$ws1.Range("F9").Value = 371
$ws1.Range("F10").Value = 4
$ws1.Range("F11").Value = 488
$ws1.Range("F13").Value = 150
$ws1.Range("F14").Value = 12279
$ws1.Range("F15").Value = 66
$ws1.Range("F16").Value = 5463

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 112

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 112
$ws4.Range("F7").Value = 1025
$ws4.Range("F8").Value = 164
$ws4.Range("F10").Value = 194
$ws4.Range("F11").Value = 371
$ws4.Range("F12").Value = 4
$ws4.Range("F13").Value = 488
$ws4.Range("F15").Value = 150
$ws4.Range("F16").Value = 12279
$ws4.Range("F18").Value = 66
$ws4.Range("F19").Value = 5463
